# Generate Report for Handback
# Update generated timestamps (and one status value ht -> mt) to reflect
# the latest handback/generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview.Range("G2").Value = "2016-08-12 02:38:29"
$wsOverview.Range("G4").Value = "2016-08-12 02:38:29"

# --- zh-cn sheet ---
# Status column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-08-12 02:38:23"
$wsZhCn.Range("H4").Value = "2016-08-12 02:38:23"

# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-08-12 02:38:40"
$wsZhCn.Range("K4").Value = "2016-08-12 02:38:40"

# --- de-de sheet ---
# Status column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-08-12 02:38:48"
$wsDeDe.Range("K4").Value = "2016-08-12 02:38:48"
